$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46078

$ws.Range("C3").Value = 46078

$ws.Range("C4").Value = 46078

$ws.Range("A5").Value = "A 1379-2024"
$ws.Range("B5").Value = 45303.55193287037
$ws.Range("C5").Value = 46078
$ws.Range("F5").Value = "Kommuner"
$ws.Range("G5").Value = 1.5

$ws.Range("A6").Value = "A 61627-2024"
$ws.Range("B6").Value = 45646.66263888889
$ws.Range("C6").Value = 46078
$ws.Range("F6").Value = "Kommuner"
$ws.Range("G6").Value = 2.1

$ws.Range("A7").Value = "A 23295-2025"
$ws.Range("B7").Value = 45791.59071759259
$ws.Range("C7").Value = 46078
$ws.Range("G7").Value = 2.2

$ws.Range("A8").Value = "A 23294-2025"
$ws.Range("B8").Value = 45791.58967592593
$ws.Range("C8").Value = 46078
$ws.Range("G8").Value = 4.6

$ws.Range("A9").Value = "A 23301-2025"
$ws.Range("B9").Value = 45791.59498842592
$ws.Range("C9").Value = 46078
$ws.Range("G9").Value = 2.9

$ws.Range("A10").Value = "A 42462-2025"
$ws.Range("B10").Value = 45905.45730324074
$ws.Range("C10").Value = 46078
$ws.Range("G10").Value = 2.3

$ws.Range("A11").Value = "A 45167-2025"
$ws.Range("B11").Value = 45919.49364583333
$ws.Range("C11").Value = 46078
$ws.Range("G11").Value = 0.6

$ws.Range("A12").Value = "A 45088-2025"
$ws.Range("B12").Value = 45919.37598379629
$ws.Range("C12").Value = 46078
$ws.Range("G12").Value = 0.9

$ws.Range("A13").Value = "A 45158-2025"
$ws.Range("B13").Value = 45919.48245370371
$ws.Range("C13").Value = 46078
$ws.Range("G13").Value = 2.7

$ws.Range("A14").Value = "A 1103-2025"
$ws.Range("B14").Value = 45666
$ws.Range("C14").Value = 46078
$ws.Range("G14").Value = 1.6

$ws.Range("A15").Value = "A 21141-2023"
$ws.Range("B15").Value = 45062
$ws.Range("C15").Value = 46078
$ws.Range("G15").Value = 3.4

$ws.Range("A16").Value = "A 34508-2025"
$ws.Range("B16").Value = 45847.44315972222
$ws.Range("C16").Value = 46078
$ws.Range("G16").Value = 2

$ws.Range("A17").Value = "A 34394-2025"
$ws.Range("B17").Value = 45846.58854166666
$ws.Range("C17").Value = 46078
$ws.Range("G17").Value = 3.9

$ws.Range("A18").Value = "A 32577-2025"
$ws.Range("B18").Value = 45838
$ws.Range("C18").Value = 46078
$ws.Range("G18").Value = 5.3

$ws.Range("A19").Value = "A 34939-2025"
$ws.Range("B19").Value = 45849.58229166667
$ws.Range("C19").Value = 46078
$ws.Range("G19").Value = 7.8

$ws.Range("A20").Value = "A 34963-2025"
$ws.Range("B20").Value = 45849.63219907408
$ws.Range("C20").Value = 46078
$ws.Range("G20").Value = 1.1

$ws.Range("A21").Value = "A 34591-2022"
$ws.Range("B21").Value = 44795.3778587963
$ws.Range("C21").Value = 46078
$ws.Range("G21").Value = 2.5

$ws.Range("A22").Value = "A 60718-2022"
$ws.Range("B22").Value = 44912.89109953704
$ws.Range("C22").Value = 46078
$ws.Range("G22").Value = 1

$ws.Range("A23").Value = "A 21379-2023"
$ws.Range("B23").Value = 45063.34819444444
$ws.Range("C23").Value = 46078
$ws.Range("F23").Value = ""
$ws.Range("G23").Value = 5.8

$ws.Range("A24").Value = "A 57000-2025"
$ws.Range("B24").Value = 45977
$ws.Range("C24").Value = 46078
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = 2.3

$ws.Range("A25").Value = "A 57655-2025"
$ws.Range("B25").Value = 45981.40369212963
$ws.Range("C25").Value = 46078
$ws.Range("G25").Value = 3.4

$ws.Range("A26").Value = "A 2434-2026"
$ws.Range("B26").Value = 46036.86722222222
$ws.Range("C26").Value = 46078
$ws.Range("G26").Value = 1.2

$ws.Range("A27").Value = "A 2433-2026"
$ws.Range("B27").Value = 46036.86631944445
$ws.Range("C27").Value = 46078
$ws.Range("G27").Value = 2.6

$ws.Range("A28").Value = "A 60392-2022"
$ws.Range("B28").Value = 44910
$ws.Range("C28").Value = 46078
$ws.Range("G28").Value = 4.1

$ws.Range("A29").Value = "A 7787-2023"
$ws.Range("B29").Value = 44973
$ws.Range("C29").Value = 46078
$ws.Range("G29").Value = 5.7

$ws.Range("A30").Value = "A 60717-2022"
$ws.Range("B30").Value = 44912.89078703704
$ws.Range("C30").Value = 46078
$ws.Range("G30").Value = 0.7

$ws.Range("A31").Value = "A 43229-2024"
$ws.Range("B31").Value = 45567.88697916667
$ws.Range("C31").Value = 46078
$ws.Range("G31").Value = 9.699999999999999

$ws.Range("A32").Value = "A 59011-2025"
$ws.Range("B32").Value = 45987
$ws.Range("C32").Value = 46078
$ws.Range("G32").Value = 2.7

$ws.Range("A33").Value = "A 59192-2022"
$ws.Range("B33").Value = 44896
$ws.Range("C33").Value = 46078
$ws.Range("G33").Value = 1.5

$ws.Range("A34").Value = "A 7791-2023"
$ws.Range("B34").Value = 44973
$ws.Range("C34").Value = 46078
$ws.Range("G34").Value = 3.1

$ws.Range("A35").Value = "A 769-2023"
$ws.Range("B35").Value = 44931
$ws.Range("C35").Value = 46078
$ws.Range("G35").Value = 1.7

$ws.Range("A36").Value = "A 1486-2023"
$ws.Range("B36").Value = 44937
$ws.Range("C36").Value = 46078
$ws.Range("G36").Value = 0.2

$ws.Range("A37").Value = "A 60793-2023"
$ws.Range("B37").Value = 45260.6534375
$ws.Range("C37").Value = 46078
$ws.Range("G37").Value = 0.7

$ws.Range("A38").Value = "A 10340-2026"
$ws.Range("C38").Value = 46078
$ws.Range("G38").Value = 1

$ws.Range("A39").Value = "A 10354-2026"
$ws.Range("B39").Value = 46075
$ws.Range("C39").Value = 46078
$ws.Range("G39").Value = 0.7

$ws.Range("A40").Value = "A 57410-2022"
$ws.Range("B40").Value = 44896
$ws.Range("C40").Value = 46078
$ws.Range("G40").Value = 7.5

$ws.Range("A41").Value = "A 31486-2021"
$ws.Range("B41").Value = 44369.43783564815
$ws.Range("C41").Value = 46078
$ws.Range("G41").Value = 4.8

$ws.Range("A42").Value = "A 63664-2023"
$ws.Range("B42").Value = 45275.62074074074
$ws.Range("C42").Value = 46078
$ws.Range("G42").Value = 3.2

$ws.Range("A43").Value = "A 57955-2024"
$ws.Range("B43").Value = 45631.56939814815
$ws.Range("C43").Value = 46078
$ws.Range("G43").Value = 0.9
